$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Harrison X.%Bai%NULL%0, Ben%Hsieh%NULL%2, Zeng%Xiong%NULL%4, Zeng%Xiong%NULL%0, Kasey%Halsey%NULL%4, Kasey%Halsey%NULL%0, Ji Whae%Choi%NULL%4, Ji Whae%Choi%NULL%0, Thi My Linh%Tran%NULL%4, Thi My Linh%Tran%NULL%0, Ian%Pan%NULL%4, Ian%Pan%NULL%0, Lin-Bo%Shi%NULL%4, Lin-Bo%Shi%NULL%0, Dong-Cui%Wang%NULL%2, Ji%Mei%NULL%4, Ji%Mei%NULL%0, Xiao-Long%Jiang%NULL%2, Qiu-Hua%Zeng%NULL%2, Thomas K.%Egglin%NULL%2, Ping-Feng%Hu%NULL%4, Ping-Feng%Hu%NULL%0, Saurabh%Agarwal%NULL%2, Fangfang%Xie%NULL%4, Fangfang%Xie%NULL%0, Sha%Li%NULL%2, Terrance%Healey%NULL%4, Terrance%Healey%NULL%0, Michael K.%Atalay%NULL%2, Wei-Hua%Liao%liaoweihua2017@163.com%4, Wei-Hua%Liao%liaoweihua2017@163.com%0]'
$ws.Range("I2").Value = ''
$ws.Range("J2").Value = 'Radiological Society of North America'
$ws.Range("E3").Value = '[Tao%Ai%NULL%0, Zhenlu%Yang%NULL%0, Zhenlu%Yang%NULL%0, Hongyan%Hou%NULL%0, Hongyan%Hou%NULL%0, Chenao%Zhan%NULL%0, Chong%Chen%NULL%0, Wenzhi%Lv%NULL%0, Wenzhi%Lv%NULL%0, Qian%Tao%NULL%0, Qian%Tao%NULL%0, Ziyong%Sun%NULL%0, Ziyong%Sun%NULL%0, Liming%Xia%xialiming2017@outlook.com%0]'
$ws.Range("I3").Value = ''
$ws.Range("J3").Value = 'Radiological Society of North America'
$ws.Range("E4").Value = '[Damiano%Caruso%NULL%0, Marta%Zerunian%NULL%4, Marta%Zerunian%NULL%0, Michela%Polici%NULL%4, Michela%Polici%NULL%0, Francesco%Pucciarelli%NULL%4, Francesco%Pucciarelli%NULL%0, Tiziano%Polidori%NULL%4, Tiziano%Polidori%NULL%0, Carlotta%Rucci%NULL%4, Carlotta%Rucci%NULL%0, Gisella%Guido%NULL%4, Gisella%Guido%NULL%0, Benedetta%Bracci%NULL%4, Benedetta%Bracci%NULL%0, Chiara%de Dominicis%NULL%4, Chiara%de Dominicis%NULL%0, Andrea%Laghi%andrea.laghi@uniroma1.it%4, Andrea%Laghi%andrea.laghi@uniroma1.it%0]'
$ws.Range("I4").Value = ''
$ws.Range("J4").Value = 'Radiological Society of North America'
$ws.Range("E5").Value = '[Xiaofeng%Chen%NULL%0, Yanyan%Tang%NULL%2, Yongkang%Mo%NULL%2, Shengkai%Li%NULL%2, Daiying%Lin%NULL%2, Zhijian%Yang%NULL%2, Zhiqi%Yang%NULL%2, Hongfu%Sun%NULL%2, Jinming%Qiu%NULL%2, Yuting%Liao%NULL%2, Jianning%Xiao%NULL%2, Xiangguang%Chen%NULL%2, Xianheng%Wu%NULL%2, Renhua%Wu%NULL%2, Zhuozhi%Dai%zhuozhi@ualberta.ca%2]'
$ws.Range("I5").Value = ''
$ws.Range("J5").Value = 'Springer Berlin Heidelberg'
$ws.Range("C6").Value = 'Unknown Title'
$ws.Range("E6").Value = '[]'
$ws.Range("F6").Value = 'not found'
$ws.Range("G6").Value = 'N/A'
$ws.Range("H6").Value = "'1970-01-01"
$ws.Range("H6").ClearFormats()
$ws.Range("J6").Value = ''
$ws.Range("E7").Value = '[Hyewon%Choi%NULL%0, Xiaolong%Qi%NULL%5, Xiaolong%Qi%NULL%0, Soon Ho%Yoon%yshoka@gmail.com%4, Soon Ho%Yoon%yshoka@gmail.com%0, Sang Joon%Park%NULL%4, Sang Joon%Park%NULL%0, Kyung Hee%Lee%NULL%4, Kyung Hee%Lee%NULL%0, Jin Yong%Kim%NULL%0, Jin Yong%Kim%NULL%0, Young Kyung%Lee%NULL%4, Young Kyung%Lee%NULL%0, Hongseok%Ko%NULL%4, Hongseok%Ko%NULL%0, Ki Hwan%Kim%NULL%4, Ki Hwan%Kim%NULL%0, Chang Min%Park%NULL%4, Chang Min%Park%NULL%0, Yun-Hyeon%Kim%NULL%4, Yun-Hyeon%Kim%NULL%0, Junqiang%Lei%NULL%0, Junqiang%Lei%NULL%0, Jung Hee%Hong%NULL%4, Jung Hee%Hong%NULL%0, Hyungjin%Kim%NULL%4, Hyungjin%Kim%NULL%0, Eui Jin%Hwang%NULL%4, Eui Jin%Hwang%NULL%0, Seung Jin%Yoo%NULL%4, Seung Jin%Yoo%NULL%0, Ju Gang%Nam%NULL%4, Ju Gang%Nam%NULL%0, Chang Hyun%Lee%NULL%4, Chang Hyun%Lee%NULL%0, Jin Mo%Goo%NULL%4, Jin Mo%Goo%NULL%0]'
$ws.Range("I7").Value = ''
$ws.Range("J7").Value = 'Radiological Society of North America'
$ws.Range("E8").Value = '[Yuki%Himoto%yukihimoto@gmail.com%0, Akihiko%Sakata%NULL%2, Mitsuhiro%Kirita%NULL%2, Takashi%Hiroi%NULL%2, Ken-ichiro%Kobayashi%NULL%2, Kenji%Kubo%NULL%2, Hyunjin%Kim%NULL%2, Azusa%Nishimoto%NULL%2, Chikara%Maeda%NULL%2, Akira%Kawamura%NULL%2, Nobuhiro%Komiya%NULL%2, Shigeaki%Umeoka%NULL%2]'
$ws.Range("I8").Value = ''
$ws.Range("J8").Value = 'Springer Singapore'
$ws.Range("E9").Value = '[Chunqin%Long%NULL%0, Huaxiang%Xu%NULL%2, Qinglin%Shen%NULL%2, Xianghai%Zhang%NULL%2, Bing%Fan%26171381@qq.com%3, Chuanhong%Wang%NULL%2, Bingliang%Zeng%NULL%2, Zicong%Li%NULL%2, Xiaofen%Li%NULL%3, Honglu%Li%NULL%2]'
$ws.Range("I9").Value = ''
$ws.Range("J9").Value = 'Elsevier B.V.'
$ws.Range("E10").Value = '[Congliang%Miao%NULL%0, Mengdi%Jin%NULL%2, Li%Miao%NULL%2, Xinying%Yang%NULL%2, Peng%Huang%NULL%3, Huanwen%Xiong%NULL%2, Peijie%Huang%NULL%2, Qi%Zhao%NULL%2, Jiang%Du%NULL%0, Jiang%Hong%NULL%2]'
$ws.Range("I10").Value = ''
$ws.Range("J10").Value = 'Elsevier Inc.'
$ws.Range("E11").Value = '[Zeying%Wen%NULL%0, Yonge%Chi%NULL%4, Yonge%Chi%NULL%0, Liang%Zhang%NULL%4, Liang%Zhang%NULL%0, Huan%Liu%NULL%5, Huan%Liu%NULL%0, Kun%Du%NULL%2, Zhengxing%Li%NULL%4, Zhengxing%Li%NULL%0, Jie%Chen%NULL%0, Jie%Chen%NULL%0, Liuhui%Cheng%NULL%4, Liuhui%Cheng%NULL%0, Daoqing%Wang%wangdaoqing1215@126.com%4, Daoqing%Wang%wangdaoqing1215@126.com%0]'
$ws.Range("I11").Value = ''
$ws.Range("J11").Value = 'Radiological Society of North America'
$ws.Range("E12").Value = '[Wanbo%Zhu%NULL%0, Kai%Xie%NULL%0, Kai%Xie%NULL%0, Hui%Lu%NULL%0, Lei%Xu%bayinhexl@126.com%0, Shusheng%Zhou%zhouss108@163.com%0, Shiyuan%Fang%fangshiyuan2008@126.com%0]'
$ws.Range("I12").Value = ''
$ws.Range("J12").Value = 'John Wiley and Sons Inc.'
$ws.Range("E13").Value = '[Anthony%Dangis%NULL%0, Christopher%Gieraerts%NULL%4, Christopher%Gieraerts%NULL%0, Yves%De Bruecker%NULL%4, Yves%De Bruecker%NULL%0, Lode%Janssen%NULL%4, Lode%Janssen%NULL%0, Hanne%Valgaeren%NULL%4, Hanne%Valgaeren%NULL%0, Dagmar%Obbels%NULL%4, Dagmar%Obbels%NULL%0, Marc%Gillis%NULL%4, Marc%Gillis%NULL%0, Marc%Van Ranst%NULL%0, Marc%Van Ranst%NULL%0, Johan%Frans%NULL%3, Annick%Demeyere%NULL%4, Annick%Demeyere%NULL%0, Rolf%Symons%rolf.symons@imelda.be%5, Rolf%Symons%rolf.symons@imelda.be%0]'
$ws.Range("I13").Value = ''
$ws.Range("J13").Value = 'Radiological Society of North America'
